# Add a new "circuit" column (U) with Illinois judicial circuit names for each
# courthouse row, and drop a stray fill-style override that had been applied
# to Q8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column U
$ws.Cells.Item(1, 21).Value = "circuit"

# Circuit name for each data row (rows 2-111, i.e. courthouses 1-110)
$circuits = @(
    "Eighth",
    "First",
    "Third",
    "Seventeenth",
    "Eighth",
    "Thirteenth",
    "Eighth",
    "Fifthteenth",
    "Eighth",
    "Sixth",
    "Fourth",
    "Fifth",
    "Fourth",
    "Fourth",
    "Fifith",
    "Cook",
    "Cook",
    "Cook",
    "Cook",
    "Cook",
    "Cook",
    "Cook",
    "Cook",
    "Cook",
    "Second",
    "Fifth",
    "Twenty-third",
    "Sixth",
    "Sixth",
    "Eighteenth",
    "Fifth",
    "Second",
    "Fourth",
    "Fourth",
    "Eleventh",
    "Second",
    "Ninth",
    "Second",
    "Seventh",
    "Thirteenth",
    "Second",
    "Ninth",
    "Second",
    "Ninth",
    "Fourteenth",
    "Twenty-first",
    "First",
    "Fourth",
    "Second",
    "Seventh",
    "Fifthteenth",
    "First",
    "Sixteenth",
    "Twenty-first",
    "Twenty-third",
    "Ninth",
    "Nineteenth",
    "Thirteenth",
    "Second",
    "Fifteenth",
    "Eleventh",
    "Eleventh",
    "Sixth",
    "Seventh",
    "Third",
    "Fourth",
    "Tenth",
    "Eighth",
    "First",
    "Ninth",
    "Twenty-second",
    "Eleventh",
    "Eighth",
    "Fourteenth",
    "Twenty-fourth",
    "Fourth",
    "Seventh",
    "Sixth",
    "Fifteenth",
    "Tenth",
    "Twenty-fourth",
    "Sixth",
    "Eighth",
    "First",
    "First",
    "Tenth",
    "Twenty-fourth",
    "Second",
    "Fourteenth",
    "First",
    "Seventh",
    "Eighth",
    "Seventh",
    "Fourth",
    "Twentieth",
    "Tenth",
    "Fifteenth",
    "Tenth",
    "First",
    "Fifth",
    "Second",
    "Ninth",
    "Twenty-fourth",
    "Second",
    "Second",
    "Fourteenth",
    "Twelfth",
    "First",
    "Seventeenth",
    "Eleventh"
)

for ($i = 0; $i -lt $circuits.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 21).Value = $circuits[$i]
}

# Remove the accidental "applyFill" style that had been left on Q8
$ws.Range("Q8").ClearFormats()
